$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format first so numeric-looking strings
# (e.g. "310.23", "0.4682") are stored as text, matching the source data
# (inline/shared strings), not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.946.54"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").Value = "1.817.72"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "310.23"
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").Value = "0.4682"
$ws.Range("D8").Value = "0.3666"
$ws.Range("E8").Value = "  -0.68%  "
$ws.Range("D9").Value = "0.07350"
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").Value = "0.8720"
$ws.Range("E10").Value = "  +0.11%  "
$ws.Range("D11").Value = "20.26"
$ws.Range("E11").Value = "  -0.91%  "
$ws.Range("D12").Value = "1.824.13"
$ws.Range("E12").Value = "  -2.20%  "
$ws.Range("D13").Value = "5.403"
$ws.Range("E13").Value = "  +0.97%  "
$ws.Range("D14").Value = "0.07109"
$ws.Range("E14").Value = "  +0.82%  "
$ws.Range("D15").Value = "6.512"
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("E16").Value = "  +0.16%  "
$ws.Range("E17").Value = "  +0.10%  "
$ws.Range("D18").Value = "0.000008706"
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("D20").Value = "14.64"
$ws.Range("E20").Value = "  -0.28%  "
$ws.Range("D21").Value = "26.965.93"
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("D22").Value = "5.285"
$ws.Range("E22").Value = "  -0.67%  "
$ws.Range("E23").Value = "  +0.83%  "
$ws.Range("D24").Value = "2.040.30"
$ws.Range("E24").Value = "  -1.92%  "
$ws.Range("D25").Value = "1.894"
$ws.Range("E25").Value = "  -0.61%  "
$ws.Range("D26").Value = "150.87"
$ws.Range("E26").Value = "  -0.48%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("E28").Value = "  +0.72%  "
$ws.Range("E29").Value = "  -0.95%  "
$ws.Range("D30").Value = "116.95"
$ws.Range("E30").Value = "  +0.99%  "
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("D32").Value = "0.7585"
$ws.Range("E32").Value = "  +0.92%  "
$ws.Range("D33").Value = "1.161"
$ws.Range("E33").Value = "  +0.81%  "
$ws.Range("D34").Value = "4.507"
$ws.Range("E34").Value = "  +1.13%  "
$ws.Range("D35").Value = "2.910"
$ws.Range("E35").Value = "  -0.29%  "
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("D37").Value = "1.099"
$ws.Range("E37").Value = "  +0.26%  "
$ws.Range("E38").Value = "  +1.09%  "
$ws.Range("D39").Value = "0.01944"
$ws.Range("E39").Value = "  -0.68%  "
$ws.Range("D40").Value = "2.965"
$ws.Range("E40").Value = "  +1.28%  "
$ws.Range("D41").Value = "2.387"
$ws.Range("E41").Value = "  -1.36%  "
$ws.Range("D42").Value = "0.5298"
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("D43").Value = "7.158"
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").Value = "0.1653"
$ws.Range("E44").Value = "  -0.55%  "
$ws.Range("D45").Value = "8.440"
$ws.Range("E45").Value = "  +0.10%  "
$ws.Range("E46").Value = "  -1.17%  "
$ws.Range("E47").Value = "  +2.07%  "
$ws.Range("E48").Value = "  +0.16%  "
$ws.Range("D49").Value = "103.40"
$ws.Range("E49").Value = "  +0.39%  "
$ws.Range("D50").Value = "1.662"
$ws.Range("E50").Value = "  -0.32%  "
$ws.Range("D51").Value = "0.06300"
$ws.Range("E51").Value = "  +0.34%  "

# Remove the temporary number-format override so the cell style
# index matches the original (unstyled) cells.
$ws.Range("D2:D51").ClearFormats()
